$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 428494.8
$ws.Range("I2").Value = 661491.4399999999
$ws.Range("J2").Value = 1334.3334
$ws.Range("K2").Value = 661491.4399999999
$ws.Range("L2").Value = 1334.3334
$ws.Range("M2").Value = -661378.4399999999
$ws.Range("N2").Value = -1560.3334
$ws.Range("H9").Value = 350.76923
$ws.Range("I9").Value = 83.71429000000001
$ws.Range("K9").Value = 83.71429000000001
$ws.Range("M9").Value = 85.28570999999999
$ws.Range("H18").Value = 939.6
$ws.Range("I18").Value = 1049.75
$ws.Range("J18").Value = 499
$ws.Range("K18").Value = 1049.75
$ws.Range("L18").Value = 499
$ws.Range("M18").Value = -765.75
$ws.Range("N18").Value = -1067
$ws.Range("H32").Value = 3318.6
$ws.Range("I32").Value = 2799.5
$ws.Range("J32").Value = 3664.6667
$ws.Range("K32").Value = 2799.5
$ws.Range("L32").Value = 3664.6667
$ws.Range("M32").Value = -2473.5
$ws.Range("N32").Value = -4316.6667
$ws.Range("H41").Value = 321.17392
$ws.Range("I41").Value = 441.85715
$ws.Range("K41").Value = 441.85715
$ws.Range("M41").Value = -1.85714999999999
$ws.Range("H106").Value = 5374.3125
$ws.Range("I106").Value = 3386.7144
$ws.Range("K106").Value = 3386.7144
$ws.Range("M106").Value = -2755.7144
$ws.Range("H116").Value = 4036.4
$ws.Range("J116").Value = 3497.5
$ws.Range("L116").Value = 3497.5
$ws.Range("N116").Value = -10381.5
$ws.Range("H132").Value = 5614.15
$ws.Range("I132").Value = 5143.375
$ws.Range("J132").Value = 7497.25
$ws.Range("K132").Value = 15430.125
$ws.Range("L132").Value = 22491.75
$ws.Range("M132").Value = -12900.125
$ws.Range("N132").Value = -27551.75
$ws.Range("H137").Value = 28793.19
$ws.Range("I137").Value = 41366.88
$ws.Range("J137").Value = 2598
$ws.Range("K137").Value = 124100.64
$ws.Range("L137").Value = 7794
$ws.Range("M137").Value = -121550.64
$ws.Range("N137").Value = -12894
$ws.Range("H138").Value = 1234.75
$ws.Range("I138").Value = 958.5294
$ws.Range("K138").Value = 2875.5882
$ws.Range("M138").Value = 2264.4118
$ws.Range("H141").Value = 1510.84
$ws.Range("I141").Value = 1510.84
$ws.Range("K141").Value = 4532.52
$ws.Range("M141").Value = 647.4800000000005
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8773
$ws.Range("I31").Value = 8773
$ws.Range("K31").Value = 8773
$ws.Range("M31").Value = -8479
$ws.Range("H61").Value = 9709.714
$ws.Range("I61").Value = 9902.77
$ws.Range("J61").Value = 7200
$ws.Range("K61").Value = 9902.77
$ws.Range("L61").Value = 7200
$ws.Range("M61").Value = -9690.77
$ws.Range("N61").Value = -7624
$ws.Range("H132").Value = 60354.445
$ws.Range("I132").Value = 67051.69
$ws.Range("K132").Value = 201155.07
$ws.Range("M132").Value = -198625.07
$ws.Range("H136").Value = 9709.714
$ws.Range("I136").Value = 9902.77
$ws.Range("J136").Value = 7200
$ws.Range("K136").Value = 29708.31
$ws.Range("L136").Value = 21600
$ws.Range("M136").Value = -27158.31
$ws.Range("N136").Value = -26700
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3283.818
$ws.Range("I99").Value = 3818
$ws.Range("J99").Value = 1859.3334
$ws.Range("K99").Value = 3818
$ws.Range("L99").Value = 1859.3334
$ws.Range("M99").Value = -2320
$ws.Range("N99").Value = -4855.3334
$ws.Range("H105").Value = 2985.3684
$ws.Range("J105").Value = 3511.8572
$ws.Range("L105").Value = 3511.8572
$ws.Range("N105").Value = -7005.8572
$ws.Range("H134").Value = 2549.875
$ws.Range("I134").Value = 2310.5715
$ws.Range("K134").Value = 6931.7145
$ws.Range("M134").Value = -4396.7145
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 1000
$ws.Range("I51").Value = 1000
$ws.Range("K51").Value = 1000
$ws.Range("M51").Value = -264
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -652
$ws.Range("H132").Value = 1635.619
$ws.Range("I132").Value = 1657.9
$ws.Range("J132").Value = 1190
$ws.Range("K132").Value = 4973.700000000001
$ws.Range("L132").Value = 3570
$ws.Range("M132").Value = -2443.700000000001
$ws.Range("N132").Value = -8630
$ws.Range("H134").Value = 49476.523
$ws.Range("I134").Value = 64131.75
$ws.Range("K134").Value = 192395.25
$ws.Range("M134").Value = -189860.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 838210.5600000001
$ws.Range("I4").Value = 555662.8
$ws.Range("K4").Value = 1666988.4
$ws.Range("M4").Value = -1666876.4
$ws.Range("H86").Value = 507.5
$ws.Range("I86").Value = 429
$ws.Range("K86").Value = 1287
$ws.Range("M86").Value = -101
$ws.Range("H89").Value = 507.5
$ws.Range("I89").Value = 429
$ws.Range("K89").Value = 3861
$ws.Range("M89").Value = 2067
$ws.Range("H92").Value = 333.42856
$ws.Range("I92").Value = 316.1111
$ws.Range("J92").Value = 364.6
$ws.Range("K92").Value = 948.3333
$ws.Range("L92").Value = 1093.8
$ws.Range("M92").Value = 299.6667
$ws.Range("N92").Value = -3589.8
$ws.Range("H137").Value = 2096.625
$ws.Range("I137").Value = 1710.5
$ws.Range("J137").Value = 3255
$ws.Range("K137").Value = 5131.5
$ws.Range("L137").Value = 9765
$ws.Range("M137").Value = -31.5
$ws.Range("N137").Value = -19965
$ws.Range("H140").Value = 2318.8333
$ws.Range("J140").Value = 4250
$ws.Range("L140").Value = 12750
$ws.Range("N140").Value = -23110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 16556.518
$ws.Range("J24").Value = 16556.518
$ws.Range("L24").Value = 16556.518
$ws.Range("N24").Value = -16902.518
$ws.Range("H107").Value = 48766.715
$ws.Range("I107").Value = 63357.75
$ws.Range("K107").Value = 63357.75
$ws.Range("M107").Value = -61437.75
$ws.Range("H126").Value = 8181.3076
$ws.Range("J126").Value = 8395.4
$ws.Range("L126").Value = 25186.2
$ws.Range("N126").Value = -30126.2
$ws.Range("H132").Value = 89468.08
$ws.Range("I132").Value = 104189.55
$ws.Range("K132").Value = 312568.65
$ws.Range("M132").Value = -310038.65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3096.2285
$ws.Range("I16").Value = 2994.7585
$ws.Range("J16").Value = 3586.6667
$ws.Range("K16").Value = 2994.7585
$ws.Range("L16").Value = 3586.6667
$ws.Range("M16").Value = -2824.7585
$ws.Range("N16").Value = -3926.6667
$ws.Range("H22").Value = 2745.4119
$ws.Range("J22").Value = 4010.3333
$ws.Range("L22").Value = 4010.3333
$ws.Range("N22").Value = -4600.3333
$ws.Range("H27").Value = 2745.4119
$ws.Range("J27").Value = 4010.3333
$ws.Range("L27").Value = 4010.3333
$ws.Range("N27").Value = -4224.3333
$ws.Range("H40").Value = 28137.3
$ws.Range("I40").Value = 44315.4
$ws.Range("J40").Value = 11959.2
$ws.Range("K40").Value = 44315.4
$ws.Range("L40").Value = 11959.2
$ws.Range("M40").Value = -44179.4
$ws.Range("N40").Value = -12231.2
$ws.Range("H46").Value = 3897
$ws.Range("J46").Value = 5757.8335
$ws.Range("L46").Value = 5757.8335
$ws.Range("N46").Value = -6133.8335
$ws.Range("H55").Value = 870.871
$ws.Range("I55").Value = 671.619
$ws.Range("J55").Value = 1289.3
$ws.Range("K55").Value = 671.619
$ws.Range("L55").Value = 1289.3
$ws.Range("M55").Value = -498.619
$ws.Range("N55").Value = -1635.3
$ws.Range("H93").Value = 1966.3334
$ws.Range("I93").Value = 1681
$ws.Range("K93").Value = 1681
$ws.Range("M93").Value = -433
$ws.Range("H132").Value = 20519.42
$ws.Range("I132").Value = 21572.11
$ws.Range("J132").Value = 7045
$ws.Range("K132").Value = 64716.33
$ws.Range("L132").Value = 21135
$ws.Range("M132").Value = -62186.33
$ws.Range("N132").Value = -26195
$ws.Range("H136").Value = 2724.5925
$ws.Range("I136").Value = 2464
$ws.Range("K136").Value = 7392
$ws.Range("M136").Value = -4842
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 6450
$ws.Range("J18").Value = 6450
$ws.Range("L18").Value = 6450
$ws.Range("N18").Value = -6796
$ws.Range("H122").Value = 145950.28
$ws.Range("I122").Value = 2430.6
$ws.Range("K122").Value = 7291.799999999999
$ws.Range("M122").Value = -4841.799999999999
$ws.Range("H132").Value = 25666.934
$ws.Range("I132").Value = 27388.5
$ws.Range("J132").Value = 1565
$ws.Range("K132").Value = 82165.5
$ws.Range("L132").Value = 4695
$ws.Range("M132").Value = -79635.5
$ws.Range("N132").Value = -9755
